$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task description text for B12 (Abstract/Projektantrag), now continuing with
# logo/image search and template prep. Line break + wrap text like the real edit.
$ws.Range("B12").Value = "Abstract und Projektantrag erarbeiten, Logo-/Bildsuche für Arbeit" + [char]10 + "Templates vorbereiten"
$ws.Range("B12").WrapText = $true

# Extra hour worked that day.
$ws.Range("C12").Value = 9

# Date cell keeps left alignment applied explicitly (new style).
$ws.Range("A12").HorizontalAlignment = -4131

# Row now needs more height to show the wrapped text.
$ws.Rows.Item(12).RowHeight = 30

# Selection moved on to B11 for the next entry.
$ws.Range("B11").Select() | Out-Null
